$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3762.389
$ws.Range("I74").Value = 3896.6667
$ws.Range("J74").Value = 3695.25
$ws.Range("K74").Value = 3896.6667
$ws.Range("L74").Value = 3695.25
$ws.Range("M74").Value = -2960.6667
$ws.Range("N74").Value = -5567.25

$ws.Range("H76").Value = 3518.9062
$ws.Range("I76").Value = 3357.5
$ws.Range("J76").Value = 4218.3335
$ws.Range("K76").Value = 3357.5
$ws.Range("L76").Value = 4218.3335
$ws.Range("M76").Value = -3042.5
$ws.Range("N76").Value = -4848.3335

$ws.Range("H77").Value = 3762.389
$ws.Range("I77").Value = 3896.6667
$ws.Range("J77").Value = 3695.25
$ws.Range("K77").Value = 19483.3335
$ws.Range("L77").Value = 18476.25
$ws.Range("M77").Value = -14803.3335
$ws.Range("N77").Value = -27836.25

$ws.Range("H79").Value = 3518.9062
$ws.Range("I79").Value = 3357.5
$ws.Range("J79").Value = 4218.3335
$ws.Range("K79").Value = 3357.5
$ws.Range("L79").Value = 4218.3335
$ws.Range("M79").Value = -2265.5
$ws.Range("N79").Value = -6402.3335

$ws.Range("H138").Value = 1508.66
$ws.Range("I138").Value = 799.8077
$ws.Range("J138").Value = 1757.7162
$ws.Range("K138").Value = 2399.4231
$ws.Range("L138").Value = 5273.1486
$ws.Range("M138").Value = 2740.5769
$ws.Range("N138").Value = -15553.1486

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 1460.5454
$ws.Range("I63").Value = 1328.8889
$ws.Range("K63").Value = 1328.8889
$ws.Range("M63").Value = -642.8888999999999

$ws.Range("H66").Value = 1460.5454
$ws.Range("I66").Value = 1328.8889
$ws.Range("K66").Value = 6644.4445
$ws.Range("M66").Value = -3212.4445

$ws.Range("H102").Value = 2913.1538
$ws.Range("I102").Value = 2622.8572
$ws.Range("J102").Value = 3251.8333
$ws.Range("K102").Value = 2622.8572
$ws.Range("L102").Value = 3251.8333
$ws.Range("M102").Value = -1000.8572
$ws.Range("N102").Value = -6495.8333

$ws.Range("H122").Value = 3118.2942
$ws.Range("I122").Value = 3570
$ws.Range("J122").Value = 2473
$ws.Range("K122").Value = 10710
$ws.Range("L122").Value = 7419
$ws.Range("M122").Value = -8260
$ws.Range("N122").Value = -12319

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2488.3809
$ws.Range("I99").Value = 2487.2727
$ws.Range("J99").Value = 2489.6
$ws.Range("K99").Value = 2487.2727
$ws.Range("L99").Value = 2489.6
$ws.Range("M99").Value = -989.2727
$ws.Range("N99").Value = -5485.6

$ws.Range("H103").Value = 45000
$ws.Range("J103").Value = 45000
$ws.Range("L103").Value = 45000
$ws.Range("N103").Value = -47344

$ws.Range("H134").Value = 1744.74
$ws.Range("I134").Value = 1685.8937
$ws.Range("K134").Value = 5057.6811
$ws.Range("M134").Value = -2522.6811

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1530.7693
$ws.Range("I99").Value = 1600
$ws.Range("J99").Value = 1300
$ws.Range("K99").Value = 1600
$ws.Range("L99").Value = 1300
$ws.Range("M99").Value = -102
$ws.Range("N99").Value = -4296

$ws.Range("H122").Value = 14622.214
$ws.Range("I122").Value = 7068.6665
$ws.Range("J122").Value = 20287.375
$ws.Range("K122").Value = 21205.9995
$ws.Range("L122").Value = 60862.125
$ws.Range("M122").Value = -18755.9995
$ws.Range("N122").Value = -65762.125

$ws.Range("H126").Value = 1530.7693
$ws.Range("I126").Value = 1600
$ws.Range("J126").Value = 1300
$ws.Range("K126").Value = 4800
$ws.Range("L126").Value = 3900
$ws.Range("M126").Value = -2330
$ws.Range("N126").Value = -8840

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 15628973
$ws.Range("I5").Value = 659.15
$ws.Range("J5").Value = 41676164
$ws.Range("K5").Value = 1977.45
$ws.Range("L5").Value = 125028492
$ws.Range("M5").Value = -1865.45
$ws.Range("N5").Value = -125028716

$ws.Range("H135").Value = 15628973
$ws.Range("I135").Value = 659.15
$ws.Range("J135").Value = 41676164
$ws.Range("K135").Value = 5932.349999999999
$ws.Range("L135").Value = 375085476
$ws.Range("M135").Value = -3397.349999999999
$ws.Range("N135").Value = -375090546

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5796.255
$ws.Range("I70").Value = 5633.3335
$ws.Range("J70").Value = 5846.385
$ws.Range("K70").Value = 5633.3335
$ws.Range("L70").Value = 5846.385
$ws.Range("M70").Value = -5363.3335
$ws.Range("N70").Value = -6386.385

$ws.Range("H73").Value = 5796.255
$ws.Range("I73").Value = 5633.3335
$ws.Range("J73").Value = 5846.385
$ws.Range("K73").Value = 5633.3335
$ws.Range("L73").Value = 5846.385
$ws.Range("M73").Value = -4697.3335
$ws.Range("N73").Value = -7718.385

$ws.Range("H122").Value = 7741.8887
$ws.Range("J122").Value = 2096.7144
$ws.Range("L122").Value = 6290.1432
$ws.Range("N122").Value = -11190.1432

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2859
$ws.Range("I16").Value = 3098.75
$ws.Range("J16").Value = 1900
$ws.Range("K16").Value = 3098.75
$ws.Range("L16").Value = 1900
$ws.Range("M16").Value = -2928.75
$ws.Range("N16").Value = -2240

$ws.Range("H46").Value = 677
$ws.Range("I46").Value = 400
$ws.Range("J46").Value = 704.7
$ws.Range("K46").Value = 400
$ws.Range("L46").Value = 704.7
$ws.Range("M46").Value = -212
$ws.Range("N46").Value = -1080.7

$ws.Range("H68").Value = 997.8570999999999
$ws.Range("I68").Value = 999
$ws.Range("K68").Value = 999
$ws.Range("M68").Value = -250

$ws.Range("H71").Value = 997.8570999999999
$ws.Range("I71").Value = 999
$ws.Range("K71").Value = 4995
$ws.Range("M71").Value = -1251

$ws.Range("H122").Value = 4923.722
$ws.Range("I122").Value = 4606
$ws.Range("J122").Value = 5645.8184
$ws.Range("K122").Value = 13818
$ws.Range("L122").Value = 16937.4552
$ws.Range("M122").Value = -11368
$ws.Range("N122").Value = -21837.4552

$ws.Range("H132").Value = 3459.838
$ws.Range("I132").Value = 3092.4211
$ws.Range("J132").Value = 3847.6667
$ws.Range("K132").Value = 9277.263300000001
$ws.Range("L132").Value = 11543.0001
$ws.Range("M132").Value = -6747.263300000001
$ws.Range("N132").Value = -16603.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1003.9
$ws.Range("I113").Value = 382.41177
$ws.Range("J113").Value = 1463.2609
$ws.Range("K113").Value = 1147.23531
$ws.Range("L113").Value = 4389.7827
$ws.Range("M113").Value = 1022.76469
$ws.Range("N113").Value = -8729.7827

$ws.Range("H122").Value = 4605.5884
$ws.Range("I122").Value = 2980.9092
$ws.Range("J122").Value = 7584.1665
$ws.Range("K122").Value = 8942.7276
$ws.Range("L122").Value = 22752.4995
$ws.Range("M122").Value = -6492.7276
$ws.Range("N122").Value = -27652.4995

$ws.Range("H136").Value = 5038.26
$ws.Range("I136").Value = 2571.8333
$ws.Range("J136").Value = 7314.9614
$ws.Range("K136").Value = 7715.499899999999
$ws.Range("L136").Value = 21944.8842
$ws.Range("M136").Value = -5165.499899999999
$ws.Range("N136").Value = -27044.8842
